$wb = $excel.ActiveWorkbook

# Rename the "Include from ..." sheets to "Include #0" / "Include #1"
$wb.Worksheets.Item("Include from TRE_R259-HL7Part").Name = "Include #0"
$wb.Worksheets.Item("Include from TRE_R260-HL7Role").Name = "Include #1"

# Update metadata values on the Metadata sheet
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B3").Value = "2.0.0"
$ws.Range("B6").Value = "active"
$ws.Range("B8").Value = "2024-09-24T12:46:07+00:00"
